$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Update the "取得日時" (acquisition timestamp) column for all data rows (2-13)
# to the new run timestamp.
$ws.Range("A2:A13").Value = "2026-01-08 18:34:38"
